$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A16").Copy()
$ws.Range("A17").PasteSpecial(-4122)
$ws.Range("A17").Value = "bayi"

$ws.Range("A16").Copy()
$ws.Range("A18").PasteSpecial(-4122)
$ws.Range("A18").Value = "fog"
